$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 25
$ws.Range("C9").Value = 25
$ws.Range("C10").Value = 25
$ws.Range("C11").Value = 25

$ws.Range("C12").Select()
